$wb = $excel.ActiveWorkbook

# --- Step 1: rename sheets ---
# Before: LoginSheet, ContactSheet1 (sheetId 2), ContactSheet (sheetId 3)
# After:  LoginSheet, ContactSheet  (sheetId 2), ContactSheet2 (sheetId 3)
$wsOldContactSheet = $wb.Worksheets.Item("ContactSheet")
$wsOldContactSheet.Name = "ContactSheet2"

$wsMain = $wb.Worksheets.Item("ContactSheet1")
$wsMain.Name = "ContactSheet"

# --- Step 2: remember existing hyperlinks (cell + target) so we can
#     re-create them after the column shift below. We must capture the
#     data first (row/col/address) and only then delete -- deleting
#     while holding onto live Hyperlink object references confuses the
#     (re-indexing) collection, so each delete re-enumerates fresh.
$hyperlinkInfo = @()
foreach ($h in $wsMain.Hyperlinks) {
    $hyperlinkInfo += ,@($h.Range.Row, $h.Range.Column, $h.Address)
}
$hyperlinkCount = $hyperlinkInfo.Count
for ($i = 0; $i -lt $hyperlinkCount; $i++) {
    foreach ($h in $wsMain.Hyperlinks) {
        $h.Delete()
        break
    }
}

# --- Step 3: insert a new "Middle Name" column (C) with a "Sri" value ---
$bWidth = $wsMain.Columns("B").ColumnWidth
$wsMain.Columns("C").Insert()
$wsMain.Columns("C").ColumnWidth = $bWidth

$wsMain.Range("C1").Value = "Middle Name"
$wsMain.Range("C2").Value = "Sri"

# --- Step 4: re-create the hyperlinks, shifting any that were at or past
#     column C one column to the right. Adding a hyperlink re-applies the
#     "Hyperlink" cell style, so restore the original (shifted-along)
#     style afterwards to match what a plain column insert would have
#     produced.
foreach ($item in $hyperlinkInfo) {
    $row = $item[0]
    $col = $item[1]
    $addr = $item[2]
    if ($col -ge 3) { $col = $col + 1 }
    $target = $wsMain.Cells.Item($row, $col)
    $origStyle = $target.Style
    $wsMain.Hyperlinks.Add($target, $addr) | Out-Null
    $target.Style = $origStyle
}

# --- Step 5: update selection to the newly inserted cell ---
$wsMain.Range("C2").Select()
